$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# Sheet 1 - 展览
$ws1.Range("F8").Value = 1226
$ws1.Range("F11").Value = 886
$ws1.Range("F12").Value = 706
$ws1.Range("F13").Value = 192
$ws1.Range("F14").Value = 522
$ws1.Range("F18").Value = 2962
$ws1.Range("F26").Value = 5323
$ws1.Range("F28").Value = 990
$ws1.Range("F29").Value = 26
$ws1.Range("F31").Value = 332
$ws1.Range("F32").Value = 1112
$ws1.Range("F35").Value = 293

# Sheet 2 - 演出
$ws2.Range("F4").Value = 1141
$ws2.Range("F14").Value = 614
$ws2.Range("F22").Value = 41
$ws2.Range("F25").Value = 281
$ws2.Range("F26").Value = 3961

# Sheet 3 - 本地生活
$ws3.Range("F5").Value = 2480

# Sheet 4 - 全部类型
$ws4.Range("F5").Value = 2480
$ws4.Range("F15").Value = 1226
$ws4.Range("F17").Value = 886
$ws4.Range("F18").Value = 706
$ws4.Range("F19").Value = 1141
$ws4.Range("F20").Value = 1141
$ws4.Range("F21").Value = 192
$ws4.Range("F22").Value = 522
$ws4.Range("F25").Value = 2962
$ws4.Range("F32").Value = 5323
$ws4.Range("F34").Value = 990
$ws4.Range("F35").Value = 614
$ws4.Range("F36").Value = 614
$ws4.Range("F37").Value = 26
$ws4.Range("F39").Value = 332
$ws4.Range("F43").Value = 41
$ws4.Range("F46").Value = 1112
$ws4.Range("F51").Value = 293
